$d = $word.ActiveDocument

# Update the date heading
$d.Paragraphs.Item(1).Range.Text = "2025-11-13 Thursday"

# Update the 20x5 practice-problem table, row-major order
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "47-4="  # was "43-25="
$t.Cell(1, 2).Range.Text = "14+73="  # was "61-37="
$t.Cell(1, 3).Range.Text = "99-89="  # was "45-9="
$t.Cell(1, 4).Range.Text = "95-71="  # was "46+26="
$t.Cell(1, 5).Range.Text = "4+37="  # was "5+53="

$t.Cell(2, 1).Range.Text = "59+32="  # was "90-87="
$t.Cell(2, 2).Range.Text = "8+38="  # was "53-42="
$t.Cell(2, 3).Range.Text = "1+19="  # was "90-11="
$t.Cell(2, 4).Range.Text = "47+13="  # was "95-65="
$t.Cell(2, 5).Range.Text = "1+50="  # was "90-66="

$t.Cell(3, 1).Range.Text = "63+7="  # was "64+28="
$t.Cell(3, 2).Range.Text = "81-74="  # was "12+42="
$t.Cell(3, 3).Range.Text = "78-50="  # was "74-4="
$t.Cell(3, 4).Range.Text = "45+7="  # was "56-40="
$t.Cell(3, 5).Range.Text = "56-41="  # was "80-18="

$t.Cell(4, 1).Range.Text = "16+36="  # was "52-6="
$t.Cell(4, 2).Range.Text = "65-38="  # was "3+3="
$t.Cell(4, 3).Range.Text = "23-20="  # was "14+68="
$t.Cell(4, 4).Range.Text = "88-56="  # was "5+77="
$t.Cell(4, 5).Range.Text = "48-0="  # was "57-18="

$t.Cell(5, 1).Range.Text = "87-11="  # was "77-37="
$t.Cell(5, 2).Range.Text = "20+8="  # was "94-7="
$t.Cell(5, 3).Range.Text = "81-43="  # was "14-11="
$t.Cell(5, 4).Range.Text = "76-15="  # was "8+23="
$t.Cell(5, 5).Range.Text = "43+24="  # was "23+45="

$t.Cell(6, 1).Range.Text = "54+20="  # was "48-39="
$t.Cell(6, 2).Range.Text = "27-10="  # was "95-62="
$t.Cell(6, 3).Range.Text = "23+26="  # was "9+35="
$t.Cell(6, 4).Range.Text = "36-18="  # was "83-54="
$t.Cell(6, 5).Range.Text = "39-22="  # was "71-56="

$t.Cell(7, 1).Range.Text = "75+16="  # was "70+27="
$t.Cell(7, 2).Range.Text = "2+6="  # was "32-16="
$t.Cell(7, 3).Range.Text = "19+23="  # was "67+0="
$t.Cell(7, 4).Range.Text = "74-9="  # was "3+72="
$t.Cell(7, 5).Range.Text = "1+47="  # was "37-23="

$t.Cell(8, 1).Range.Text = "21+33="  # was "66+2="
$t.Cell(8, 2).Range.Text = "70-4="  # was "20+1="
$t.Cell(8, 3).Range.Text = "78-51="  # was "54+6="
$t.Cell(8, 4).Range.Text = "42+12="  # was "96-54="
$t.Cell(8, 5).Range.Text = "50-21="  # was "0+78="

$t.Cell(9, 1).Range.Text = "72-51="  # was "26+57="
$t.Cell(9, 2).Range.Text = "73+11="  # was "51-9="
$t.Cell(9, 3).Range.Text = "10+15="  # was "25+4="
$t.Cell(9, 4).Range.Text = "31+4="  # was "31+6="
$t.Cell(9, 5).Range.Text = "36-31="  # was "23+8="

$t.Cell(10, 1).Range.Text = "72-43="  # was "44+7="
$t.Cell(10, 2).Range.Text = "55-39="  # was "65-61="
$t.Cell(10, 3).Range.Text = "90-73="  # was "33+57="
$t.Cell(10, 4).Range.Text = "66-36="  # was "17+60="
$t.Cell(10, 5).Range.Text = "50+15="  # was "30+43="

$t.Cell(11, 1).Range.Text = "92-18="  # was "7+23="
$t.Cell(11, 2).Range.Text = "51-36="  # was "52+20="
$t.Cell(11, 3).Range.Text = "5+31="  # was "15+66="
$t.Cell(11, 4).Range.Text = "2+37="  # was "78+15="
$t.Cell(11, 5).Range.Text = "50+13="  # was "39-18="

$t.Cell(12, 1).Range.Text = "36+48="  # was "3+33="
$t.Cell(12, 2).Range.Text = "72-37="  # was "58-30="
$t.Cell(12, 3).Range.Text = "26-10="  # was "55-40="
$t.Cell(12, 4).Range.Text = "77-38="  # was "6+79="
$t.Cell(12, 5).Range.Text = "9-1="  # was "73-31="

$t.Cell(13, 1).Range.Text = "23+75="  # was "59-57="
$t.Cell(13, 2).Range.Text = "1+45="  # was "9+66="
$t.Cell(13, 3).Range.Text = "95-82="  # was "11+78="
$t.Cell(13, 4).Range.Text = "20+17="  # was "21-12="
$t.Cell(13, 5).Range.Text = "31+24="  # was "70-27="

$t.Cell(14, 1).Range.Text = "41+13="  # was "59-50="
$t.Cell(14, 2).Range.Text = "46-12="  # was "32+66="
$t.Cell(14, 3).Range.Text = "2+69="  # was "73-21="
$t.Cell(14, 4).Range.Text = "99-77="  # was "16+26="
$t.Cell(14, 5).Range.Text = "15+40="  # was "42+24="

$t.Cell(15, 1).Range.Text = "92-27="  # was "48+35="
$t.Cell(15, 2).Range.Text = "23+62="  # was "4+31="
$t.Cell(15, 3).Range.Text = "76-19="  # was "61+22="
$t.Cell(15, 4).Range.Text = "6+90="  # was "77-54="
$t.Cell(15, 5).Range.Text = "69-36="  # was "17+17="

$t.Cell(16, 1).Range.Text = "85-31="  # was "29+53="
$t.Cell(16, 2).Range.Text = "1+71="  # was "30+0="
$t.Cell(16, 3).Range.Text = "67+9="  # was "78-2="
$t.Cell(16, 4).Range.Text = "73-51="  # was "74-48="
$t.Cell(16, 5).Range.Text = "45+52="  # was "28+62="

$t.Cell(17, 1).Range.Text = "31+14="  # was "21-11="
$t.Cell(17, 2).Range.Text = "46-16="  # was "95-79="
$t.Cell(17, 3).Range.Text = "90-4="  # was "65-37="
$t.Cell(17, 4).Range.Text = "87+9="  # was "34+49="
$t.Cell(17, 5).Range.Text = "61+31="  # was "96-41="

$t.Cell(18, 1).Range.Text = "71-8="  # was "5+41="
$t.Cell(18, 2).Range.Text = "43-24="  # was "63-29="
$t.Cell(18, 3).Range.Text = "75+7="  # was "0+50="
$t.Cell(18, 4).Range.Text = "70+6="  # was "97-57="
$t.Cell(18, 5).Range.Text = "98-46="  # was "99-5="

$t.Cell(19, 1).Range.Text = "50+48="  # was "13-10="
$t.Cell(19, 2).Range.Text = "45+26="  # was "78-46="
$t.Cell(19, 3).Range.Text = "20+22="  # was "64-5="
$t.Cell(19, 4).Range.Text = "54+0="  # was "44+45="
$t.Cell(19, 5).Range.Text = "89-48="  # was "4+62="

$t.Cell(20, 1).Range.Text = "43+31="  # was "73-30="
$t.Cell(20, 2).Range.Text = "64-7="  # was "30-25="
$t.Cell(20, 3).Range.Text = "70-43="  # was "8+35="
$t.Cell(20, 4).Range.Text = "45+13="  # was "20+44="
$t.Cell(20, 5).Range.Text = "0+47="  # was "20+74="
